$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Preduvjeti (B5) -> "-"
$ws.Range("B5").Value = "-"
$ws.Rows(5).AutoFit()

# Posljedice - uspjesan zavrsetak (B6) -> reworded text
$ws.Range("B6").Value = "Korisniku je prikazana lista od 10 najslušanijih pjesama u posljednjem mjesecu"

# Posljedice - neuspjesan zavrsetak (B7) -> "-"
$ws.Range("B7").Value = "-"
$ws.Rows(7).AutoFit()

# Update the selection state to match the saved workbook
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("B8").Select()
